$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72") entirely; subsequent rows shift up by one.
$ws.Rows.Item(2).Delete()
